# Jogos_da_Semana_FlashScore_2024-11-20.xlsx update
# 1) Update several odds values in existing rows 3-6 (re-scraped odds).
# 2) Insert a brand-new match row at row 8 (Brazil Serie A: Bahia x Palmeiras),
#    which pushes the former rows 8-10 down to rows 9-11.
# 3) Apply additional odds updates to the rows that landed at 10 and 11
#    after the shift (these also got re-scraped odds changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Odds updates on row 3 ---
$ws.Range("S3").Value = 1.33

# --- Odds updates on row 4 ---
$ws.Range("G4").Value = 1.44
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 2
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 1.9
$ws.Range("S4").Value = 1.4
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("Y4").Value = 9
$ws.Range("AG4").Value = 15
$ws.Range("AH4").Value = 34
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 81
$ws.Range("AN4").Value = 3.25
$ws.Range("AU4").Value = 9.5
$ws.Range("AZ4").Value = 151
$ws.Range("BB4").Value = 401

# --- Odds updates on row 5 ---
$ws.Range("G5").Value = 1.65
$ws.Range("I5").Value = 5.75
$ws.Range("J5").Value = 2.25
$ws.Range("N5").Value = 9.5
$ws.Range("S5").Value = 1.4
$ws.Range("X5").Value = 7.5
$ws.Range("AG5").Value = 13
$ws.Range("AI5").Value = 17
$ws.Range("AJ5").Value = 51
$ws.Range("AL5").Value = 41
$ws.Range("AQ5").Value = 29
$ws.Range("AU5").Value = 8.5
$ws.Range("AV5").Value = 51
$ws.Range("AY5").Value = 34

# --- Odds updates on row 6 ---
$ws.Range("S6").Value = 1.53
$ws.Range("T6").Value = 2.38

# --- 2) Insert new row at 8 (shifts old rows 8-10 down to 9-11) ---
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new match data
$newRow8 = @(
    "2BqdsiBE", "20/11/2024", "18:00", "BRAZIL - SERIE A BETANO", "Bahia", "Palmeiras",
    3.6, 3.8, 1.95, 4, 2.25, 2.6, 1.04, 13, 1.25, 4, 1.75, 2.05, 1.36, 3, 1.67, 2.1,
    12, 19, 12, 41, 26, 29, 13, 7, 13, 41, 8.5, 10, 8.5, 17, 15, 23, 151,
    5.5, 19, 26, 67, 81, 151, 3, 7.5, 51, 4, 10, 19, 34, 51, 126, 501, 126
)

for ($i = 0; $i -lt $newRow8.Length; $i++) {
    $ws.Cells.Item(8, $i + 1).Value = $newRow8[$i]
}

# --- 3) Extra odds updates for the rows that shifted into positions 10 and 11 ---

# Row 10 (was old row 9: Eldense - Huesca)
$ws.Range("G10").Value = 2.3
$ws.Range("H10").Value = 2.75
$ws.Range("I10").Value = 3.8
$ws.Range("J10").Value = 3.25
$ws.Range("L10").Value = 4.75
$ws.Range("M10").Value = 1.14
$ws.Range("N10").Value = 5.5
$ws.Range("Q10").Value = 3.1
$ws.Range("R10").Value = 1.36
$ws.Range("S10").Value = 1.73
$ws.Range("T10").Value = 2
$ws.Range("X10").Value = 9
$ws.Range("Z10").Value = 21
$ws.Range("AA10").Value = 26
$ws.Range("AC10").Value = 5
$ws.Range("AG10").Value = 7.5
$ws.Range("AH10").Value = 17
$ws.Range("AW10").Value = 5.5
$ws.Range("AZ10").Value = 101

# Row 11 (was old row 10: Briton Ferry - TNS)
$ws.Range("H11").Value = 8.5
$ws.Range("J11").Value = 28
$ws.Range("K11").Value = 3.45
$ws.Range("P11").Value = 6.4
$ws.Range("Q11").Value = 1.27
$ws.Range("R11").Value = 3.4
$ws.Range("S11").Value = 1.19
$ws.Range("T11").Value = 4.2
$ws.Range("AD11").Value = 30
$ws.Range("AG11").Value = 11
$ws.Range("AH11").Value = 6.5
$ws.Range("AI11").Value = 17
$ws.Range("AK11").Value = 14
$ws.Range("AN11").Value = 37
$ws.Range("AT11").Value = 4.2
$ws.Range("AW11").Value = 3
$ws.Range("AZ11").Value = 6.6

# These four cells were (and remain) genuinely blank in row 11 (formerly row 10);
# force them back to true empty cells rather than empty-string placeholders.
$ws.Range("AQ11").ClearContents()
$ws.Range("AS11").ClearContents()
$ws.Range("BC11").ClearContents()
$ws.Range("BD11").ClearContents()

Write-Output "done"
